# Hand routing the 6 layer stackup.
# Add a new "6-layer" worksheet (copy/derivative of "8-layer"), positioned
# right after "8-layer", make it the active sheet, and populate its data.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet right after "8-layer" and rename it ---
# NOTE: worksheet handles in this runtime resolve by position, and
# Add()/Move() shuffle positions -- so fetch "8-layer" fresh (by name)
# *after* each reorder instead of reusing an older handle.
$newSheet = $wb.Worksheets.Add()
$eightLayer = $wb.Worksheets.Item("8-layer")
$newSheet.Move($null, $eightLayer)

$eightLayer = $wb.Worksheets.Item("8-layer")
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "6-layer"

# --- Column widths (best-fit-style widths, matching the target layout) ---
# (ColumnWidth is in "characters"; the engine adds ~5/6 character of padding
# when it serialises the raw OOXML <col width>, same as stock Excel COM.)
$ws.Range("A1").ColumnWidth = 3.45
$ws.Range("B1").ColumnWidth = 17.59
$ws.Range("C1").ColumnWidth = 13.88
$ws.Range("D1").ColumnWidth = 13.74
$ws.Range("E1").ColumnWidth = 3.17
$ws.Range("F1").ColumnWidth = 8.74
$ws.Range("G1").ColumnWidth = 9.45
$ws.Range("H1").ColumnWidth = 10.74
$ws.Range("I1").ColumnWidth = 11.88

# --- Header row (row 1), centered ---
$ws.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter

$ws.Range("B1").Value = "Copper Layer"
$ws.Range("C1").Value = "Dielectric Layer"
$ws.Range("D1").Value = "Material"
$ws.Range("E1").Value = "Er"
$ws.Range("F1").Value = "Thickness"
$ws.Range("G1").Value = "Desired Z0"
$ws.Range("H1").Value = "Trace Width"
$ws.Range("I1").Value = "Calculated Z0"

# --- Stackup data rows ---
# Layer 1 - top component
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "top component"
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 54

$ws.Range("C3").Value = "prepreg"
$ws.Range("D3").Value = "Grace GA-170LL"
$ws.Range("E3").Value = 4.7
$ws.Range("F3").Value = 6

# Layer 3 - ground plane
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ground plane"
$ws.Range("F4").Value = 1

$ws.Range("C5").Value = "laminate"
$ws.Range("D5").Value = "Grace GA-170LL"
$ws.Range("E5").Value = 4.7
$ws.Range("F5").Value = 6

# Layer 4 - power plane
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "power plane"
$ws.Range("F6").Value = 1

$ws.Range("C7").Value = "prepreg"
$ws.Range("D7").Value = "Grace GA-170LL"
$ws.Range("E7").Value = 4.7
$ws.Range("F7").Value = 14

# Layer 5 - power plane
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "power plane"
$ws.Range("F8").Value = 1

$ws.Range("C9").Value = "laminate"
$ws.Range("D9").Value = "Grace GA-170LL"
$ws.Range("E9").Value = 4.7
$ws.Range("F9").Value = 6

# Layer 6 - ground plane
$ws.Range("A10").Value = 6
$ws.Range("B10").Value = "ground plane"
$ws.Range("F10").Value = 1

$ws.Range("C11").Value = "prepreg"
$ws.Range("D11").Value = "Grace GA-170LL"
$ws.Range("E11").Value = 4.7
$ws.Range("F11").Value = 6

# Layer 8 - bottom component
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "bottom component"
$ws.Range("F12").Value = 1
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = 54

# --- Total thickness formula ---
$ws.Range("F19").Formula = "=+SUM(F2:F12)"

# --- Update the (previously active) "8-layer" sheet's stale selection ---
# NOTE: Range.Select() activates its parent sheet as a side effect in this
# runtime, so do this *before* activating "6-layer" (which must end up as
# the active/selected sheet).
$eightLayer.Range("A1:I23").Select()

# --- Selection / active-cell bookkeeping for the new sheet (must be last
#     so "6-layer" ends up the active tab) ---
$ws.Range("E16").Select()
$ws.Activate()

Write-Output "done"
